# Apply crypto price/volume updates per commit "Updated cryptos list on Fri Apr 12 14:44:58 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.488.50'
$ws.Range('E2').Value = '  -0.67%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.455.32'
$ws.Range('E3').Value = '  -1.41%  '

# Row 4: TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.40'
$ws.Range('E5').Value = '  +1.42%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.70'
$ws.Range('E6').Value = '  -2.44%  '

# Row 7: LidoStakedEther
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.449.47'
$ws.Range('E7').Value = '  -1.38%  '

# Row 8: XRP
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.593'
$ws.Range('E8').Value = '  -2.41%  '

# Row 9: USDC
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.02%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.191'
$ws.Range('E10').Value = '  -0.09%  '

# Row 11: Toncoin
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.05'
$ws.Range('E11').Value = '  -2.81%  '

# Row 12: Cardano
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.564'
$ws.Range('E12').Value = '  -2.74%  '

# Row 13: Avalanche
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.36'
$ws.Range('E13').Value = '  -3.87%  '

# Row 14: ShibaInu
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000269'
$ws.Range('E14').Value = '  -1.59%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.018.44'
$ws.Range('E15').Value = '  -1.18%  '

# Row 16: Polkadot
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.18'
$ws.Range('E16').Value = '  -1.02%  '

# Row 17: WrappedEther
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.462.13'
$ws.Range('E17').Value = '  -0.84%  '

# Row 18: WrappedBTC
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.614.96'

# Row 19: BitcoinCash
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '579.01'
$ws.Range('E19').Value = '  -4.14%  '

# Row 20: TRON
$ws.Range('E20').Value = '  +1.14%  '

# Row 21: Chainlink
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.17'
$ws.Range('E21').Value = '  +0.40%  '

# Row 22: Polygon
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.847'
$ws.Range('E22').Value = '  -2.22%  '

# Row 23: Uniswap
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.92'
$ws.Range('E23').Value = '  -1.96%  '

# Row 24: Litecoin
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '95.78'
$ws.Range('E24').Value = '  +0.45%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.20'
$ws.Range('E25').Value = '  -1.55%  '

# Row 26: PancakeSwap
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.62'
$ws.Range('E26').Value = '  -2.04%  '

# Row 27: Dai
$ws.Range('E27').Value = '  -0.01%  '

# Row 28: ImmutableX
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.44'
$ws.Range('E28').Value = '  -4.79%  '

# Row 29: EthereumClassic
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.81'
$ws.Range('E29').Value = '  -3.42%  '

# Row 30: RenderToken
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.62'
$ws.Range('E30').Value = '  -3.59%  '

# Row 31: Filecoin
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  -2.59%  '

# Row 32: Stacks
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.82'
$ws.Range('E32').Value = '  -6.11%  '

# Row 33: Mantle
$ws.Range('E33').Value = '  -2.46%  '

# Row 34: NEARProtocol
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.58'
$ws.Range('E34').Value = '  -4.71%  '

# Row 35: Bittensor
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '586.14'
$ws.Range('E35').Value = '  -17.19%  '

# Row 36: Cosmos
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.55'
$ws.Range('E36').Value = '  -0.96%  '

# Row 37: VeChain
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0473'
$ws.Range('E37').Value = '  +0.68%  '

# Row 38: Hedera
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0956'
$ws.Range('E38').Value = '  -3.95%  '

# Row 39: FirstDigitalUSD
$ws.Range('E39').Value = '  +0.51%  '

# Row 40: OKB
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.06'
$ws.Range('E40').Value = '  -0.66%  '

# Row 41: Kaspa
$ws.Range('E41').Value = '  -0.75%  '

# Row 42: dogwifhat
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.15'
$ws.Range('E42').Value = '  -10.87%  '

# Row 43: Maker
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.248.57'
$ws.Range('E43').Value = '  -2.14%  '

# Row 44: PEPE
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0695'
$ws.Range('E44').Value = '  +1.28%  '

# Row 45: TheGraph
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.296'
$ws.Range('E45').Value = '  -5.01%  '

# Row 46: InjectiveProtocol
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '30.92'
$ws.Range('E46').Value = '  -3.82%  '

# Row 47: ThetaToken
$ws.Range('E47').Value = '  -4.86%  '

# Row 48: Fetch.AI
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('E48').Value = '  -5.54%  '

# Row 49: Stellar
$ws.Range('E49').Value = '  -2.46%  '

# Row 50: Monero
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.75'
$ws.Range('E50').Value = '  +0.39%  '

# Row 51: USDe
$ws.Range('E51').Value = '  -0.02%  '
